# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" change
# (column E) figures with the latest values from the Fri Jul 12 2024
# GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that look like plain numbers ("530.25", "0.999", ...) are
# written with a leading apostrophe so Excel keeps storing them as text
# (as in the original file) instead of auto-converting them to numbers.

$ws.Range("D2").Value = "57.523.61"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "3.104.02"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'530.25"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'137.88"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D8").Value = "3.101.60"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.469"
$ws.Range("E9").Value = "  +4.73%  "
$ws.Range("D10").Value = "'7.31"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D12").Value = "'0.413"
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "3.631.58"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'25.62"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "57.640.01"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "3.095.75"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'6.03"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").Value = "'12.63"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "'8.07"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").Value = "'360.04"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'68.89"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "'0.504"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "0.0₃0866"
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("D29").Value = "'7.31"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").Value = "'1.87"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'6.06"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "'21.33"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "'5.08"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "'159.10"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").Value = "'6.04"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "'25.48"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "'1.27"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "'0.0670"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "2.480.41"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").Value = "'3.99"
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "'37.54"
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "3.140.95"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").Value = "'0.985"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("D49").Value = "'6.06"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "'19.73"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "'0.738"
$ws.Range("E51").Value = "  -2.92%  "
